$d = $word.ActiveDocument

# --- Part 1: remove the _GoBack bookmark that currently sits right after
# "...of the web-portal." (it moves to the end of the document, see Part 2). ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Part 2: reorder the bibliography entries so "Abril ..." comes first,
# add a new "Allen, G. H., ..." reference paragraph, keep "Raymond ..." as
# the last paragraph of the doc, drop the trailing empty paragraphs, and
# re-create the _GoBack bookmark around the (now last) "Raymond ..." text. ---

# Find the index of a paragraph whose text starts with the given prefix.
# (Re-scanning by index -- rather than caching stale Paragraph refs / using
# .Previous/.Next across mutations -- is what keeps this reliable.)
function Find-ParaIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

$raymondIdx = Find-ParaIndex "Raymond, P. A."
$abrilIdx = Find-ParaIndex "Abril, G."

# Move the "Abril ..." paragraph so that it comes right before "Raymond ...".
$insertPos = $d.Paragraphs.Item($raymondIdx).Range.Start
$d.Paragraphs.Item($abrilIdx).Range.Cut()
$target = $d.Range($insertPos, $insertPos)
$target.Paste()

# "Raymond ..." paragraph moved down by one -- re-resolve it.
$raymondIdx = Find-ParaIndex "Raymond, P. A."

# Insert a new empty paragraph right before "Raymond ..." and fill it with
# the new "Allen ..." reference.
$d.Paragraphs.Item($raymondIdx).Range.InsertParagraphBefore()
$raymondIdx = Find-ParaIndex "Raymond, P. A."
$allenPara = $d.Paragraphs.Item($raymondIdx - 1)
$allenPara.Range.Text = "Allen, G. H., and T. M. Pavelsky (2018), Global extent of rivers and streams, Science. "

# Re-resolve "Raymond ..." once more and drop everything that follows it
# (the stray blank paragraphs at the tail of the document), leaving
# "Raymond ..." as the final paragraph.
$raymondIdx = Find-ParaIndex "Raymond, P. A."
$raymondPara = $d.Paragraphs.Item($raymondIdx)
$docEnd = $d.Content.End
$tail = $d.Range($raymondPara.Range.End, $docEnd)
if ($tail.Start -lt $tail.End) {
    $tail.Delete()
}

# Re-create the _GoBack bookmark around the "Raymond ..." paragraph, which is
# now the last paragraph in the document.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
